$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows to append: date (Excel serial number), epidemiological_week,
# last_available_confirmed, last_available_deaths, new_confirmed, new_deaths
$newRows = @(
    @(44709, 0, 327439, 6348, 10, 1),
    @(44710, 0, 327442, 6348, 3, 0),
    @(44711, 0, 327451, 6348, 9, 0)
)

$dateFormat = "yyyy\-mm\-dd;@"

$startRow = 63
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = $row[0]
    $dateCell.NumberFormat = $dateFormat

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

$ws.Range("C63").Select()
